$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.7501982553528945
$ws.Range("C2").Value = 0.832014072119613
$ws.Range("D2").Value = 0.7889908256880733
$ws.Range("E2").Value = 1137

# Row 3
$ws.Range("B3").Value = 0.8539325842696629
$ws.Range("C3").Value = 0.8299531981279251
$ws.Range("D3").Value = 0.8417721518987341
$ws.Range("E3").Value = 641

# Row 4
$ws.Range("B4").Value = 0.7671957671957672
$ws.Range("C4").Value = 0.7196029776674938
$ws.Range("D4").Value = 0.7426376440460947
$ws.Range("E4").Value = 806

# Row 5
$ws.Range("B5").Value = 0.5266666666666666
$ws.Range("C5").Value = 0.4438202247191011
$ws.Range("D5").Value = 0.4817073170731707
$ws.Range("E5").Value = 356

# Row 6 (accuracy row - all columns share the same value)
$ws.Range("B6").Value = 0.7537414965986394
$ws.Range("C6").Value = 0.7537414965986394
$ws.Range("D6").Value = 0.7537414965986394
$ws.Range("E6").Value = 0.7537414965986394

# Row 7 (macro avg)
$ws.Range("B7").Value = 0.7244983183712477
$ws.Range("C7").Value = 0.7063476181585333
$ws.Range("D7").Value = 0.7137769846765181

# Row 8 (weighted avg)
$ws.Range("B8").Value = 0.7504079335191214
$ws.Range("C8").Value = 0.7537414965986394
$ws.Range("D8").Value = 0.7505824027733432
